# Daily attendance processing - 2026-01-25 04:09:33
# Normalizes the "Recorded By" (column G) attribution lists so the
# user who actually recorded the session ("System, dnasr281@gmail.com")
# is listed first, followed by the automated/system actor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Recorded By" value is the 3-name backup/system combo.
# "backup@backdoor.com, System, system" -> "system, backup@backdoor.com, System"
$threeNameRows = @(2, 28, 54)

# Rows whose "Recorded By" value is the System + dnasr281 pair.
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
$twoNameRows = @(3, 6, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 29, 32, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 55, 58, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 116, 118, 119, 120, 122, 125, 127, 135, 136, 137, 138, 142, 144, 145, 146, 148, 151, 153)

$oldThreeName = "backup@backdoor.com, System, system"
$newThreeName = "system, backup@backdoor.com, System"

$oldTwoName = "System, dnasr281@gmail.com"
$newTwoName = "dnasr281@gmail.com, System"

foreach ($r in $threeNameRows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldThreeName) {
        $cell.Value = $newThreeName
    }
}

foreach ($r in $twoNameRows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldTwoName) {
        $cell.Value = $newTwoName
    }
}
